# Commit: add a message "random object practice" in the ppt
#
# The slide that currently reads:
#   "Objectives : class practice"
# needs to become (as 4 separate runs, matching the author's diff):
#   "Objectives : " + "random " + "object and class " + "practice"
# i.e. "Objectives : random object and class practice"

$p = $ppt.ActivePresentation

$oldSentence = "Objectives : class practice"

# Locate the shape that holds the sentence, searching every slide/shape
# rather than hard-coding indices so the script is resilient to layout.
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -like "*$oldSentence*") {
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the 1-based character offset of the sentence inside the shape's text.
$charIdx0 = $tr.Text.IndexOf($oldSentence)
$startPos = $charIdx0 + 1

# New text, split exactly like the four runs in the target XML.
$seg1 = "Objectives : "
$seg2 = "random "
$seg3 = "object and class "
$seg4 = "practice"

# Replace the whole old sentence with the full new sentence first …
$whole = $tr.Characters($startPos, $oldSentence.Length)
$whole.Text = $seg1 + $seg2 + $seg3 + $seg4

# … then re-assert each segment over itself (same length, so offsets do not
# shift) which forces the text run to split at each boundary, producing the
# four separate <a:r> runs seen in the diff.
$pos = $startPos + $seg1.Length

$r2 = $tr.Characters($pos, $seg2.Length)
$r2.Text = $seg2
$pos = $pos + $seg2.Length

$r3 = $tr.Characters($pos, $seg3.Length)
$r3.Text = $seg3
$pos = $pos + $seg3.Length

$r4 = $tr.Characters($pos, $seg4.Length)
$r4.Text = $seg4
